$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: a record with an empty attachment_filename (column B)
$ws.Range("A4").Value = "0001_slr0611_right"
$ws.Range("C4").Value = "this is a row with an empty attachment file"
$ws.Range("D4").Value = "this is a row with an empty attachment file"

# New row 5: a record with an empty description (column C)
$ws.Range("B5").Value = "NC_014139.gbk"
$ws.Range("D5").Value = "This is a row with an empty description"
$ws.Range("A5").Value = "0003_slr0613_left"

# Move the active selection to A5
$ws.Range("A5").Select() | Out-Null

# Set the page to portrait orientation
$ws.PageSetup.Orientation = 1
